$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 40 (shifts rows 40:103 down to 41:104)
$ws.Rows.Item(40).Insert()

# Populate the newly inserted row 40 with the new record's data
$ws.Cells.Item(40, 1).Value = 1
$ws.Cells.Item(40, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(40, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(40, 4).Value = 45195
$ws.Cells.Item(40, 5).Value = 15
$ws.Cells.Item(40, 6).Value = 100112009
$ws.Cells.Item(40, 7).Value = "Acelga"
$ws.Cells.Item(40, 8).Value = "Sin especificar"
$ws.Cells.Item(40, 9).Value = "Primera"
$ws.Cells.Item(40, 10).Value = 250
$ws.Cells.Item(40, 11).Value = 900
$ws.Cells.Item(40, 12).Value = 1000
$ws.Cells.Item(40, 13).Value = 950
$ws.Cells.Item(40, 14).Value = "`$/atado 2,5 a 3 kilos"
$ws.Cells.Item(40, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(40, 16).Value = 317
$ws.Cells.Item(40, 17).Value = 3
$ws.Cells.Item(40, 18).Value = "Hortaliza"
